$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 26, shifting existing rows 26-41 down to 27-42
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the weekly record
$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(26, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44489
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112022
$ws.Cells.Item(26, 7).Value = "Arveja Verde"
$ws.Cells.Item(26, 8).Value = "Perfection"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 400
$ws.Cells.Item(26, 11).Value = 18000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 19000
$ws.Cells.Item(26, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 760
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Match the date style used by the other rows in column D
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
